$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.812.03'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.646.10'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.59%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.80'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("E7").Value = '  +0.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0628'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.19'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0842'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").Value = '1.870.57'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '1.645.61'
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.17'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.45'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.83%  '
$ws.Range("D17").Value = '26.825.52'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("E18").Value = '  -2.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.81'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.78%  '
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("E22").Value = '  +11.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.37'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.09'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.09'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.68'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0511'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  -2.69%  '
$ws.Range("E33").Value = '  -2.15%  '
$ws.Range("D34").Value = '1.288.50'
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.53'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("E37").Value = '  -5.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.539'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.827'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.809'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.36'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").Value = '1.797.34'
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.48'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.55'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.67'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0978'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.19%  '
